# ---------------------------------------------------------------------------
# Commit: "Wed, Jun 24, 2020  8:06:01 AM"
#
# 1) The table on slide 5 gets switched from the custom "Table_0" table
#    style over to PowerPoint's built-in "No Style, No Grid" table style.
# 2) The deck's theme colour scheme (the one and only Design/Slide Master
#    the file has) is recoloured from the "Integral / Red Violet" palette
#    back to the stock "Office" palette.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Retarget the table's style -----------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{EACEF73F-1339-4F97-AC5B-CC02406ED9F4}")

# --- 2) Recolour the theme back to the default "Office" colour scheme ------
# msoThemeColorDark1=1, Light1=2, Dark2=3, Light2=4,
# Accent1..6=5..10, Hyperlink=11, FollowedHyperlink=12
function RgbVal([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$officeColors = @{
    1  = (RgbVal 0x00 0x00 0x00)   # dk1
    2  = (RgbVal 0xFF 0xFF 0xFF)   # lt1
    3  = (RgbVal 0x44 0x54 0x6A)   # dk2
    4  = (RgbVal 0xE7 0xE6 0xE6)   # lt2
    5  = (RgbVal 0x5B 0x9B 0xD5)   # accent1
    6  = (RgbVal 0xED 0x7D 0x31)   # accent2
    7  = (RgbVal 0xA5 0xA5 0xA5)   # accent3
    8  = (RgbVal 0xFF 0xC0 0x00)   # accent4
    9  = (RgbVal 0x44 0x72 0xC4)   # accent5
    10 = (RgbVal 0x70 0xAD 0x47)   # accent6
    11 = (RgbVal 0x05 0x63 0xC1)   # hyperlink
    12 = (RgbVal 0x95 0x4F 0x72)   # followed hyperlink
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i]
}
